$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "LIVE CAM List" entries being appended to the bottom of the table
# (Monterey Bay Aquarium cams + a Chile traffic cam), mirroring the
# existing row layout: A=Category, B=lat/long, C=Location, D=CITY,
# E=COUNTRY, F=YouTube Link.
$rows = @(
    @{ A = "LIVE, AQUARIUM";      B = "36.61818880625526, -121.90155004951903"; C = "Live Sea Otter Cam - Monterey Bay Aquarium";   D = "CA";          E = "USA";   F = "abbR-Ttd-cA" },
    @{ A = "LIVE, AQUARIUM";      B = "36.61853302540104, -121.90200628992626"; C = "Live Kelp Forest Cam - Monterey Bay Aquarium"; D = "CA";          E = "USA";   F = "w3LjpFhySTg" },
    @{ A = "LIVE, AQUARIUM";      B = "36.61750418269601, -121.90111105097645"; C = "Live Open Sea Cam - Monterey Bay Aquarium";    D = "CA";          E = "USA";   F = "zM2d_2k9Nj0" },
    @{ A = "LIVE, CITY, TRAFFIC"; B = "-36.821446551229776, -73.046832180092";  C = "Concepción Centro - Ultra HD LIVE";            D = "Concepción";  E = "Chile"; F = "Z5ROgEcAWVA" }
)

$startRow = 83
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F

    # Columns A, C, D, E carry the thin left/right "grid" border used by
    # the rest of the body rows (column B/F stay borderless, matching the
    # existing sheet pattern).
    foreach ($col in @("A", "C", "D", "E")) {
        $cell = $ws.Range($col + $r)
        $cell.Borders.Item(7).Color = 0
        $cell.Borders.Item(7).Weight = 2
        $cell.Borders.Item(7).LineStyle = 1
        $cell.Borders.Item(10).Color = 0
        $cell.Borders.Item(10).Weight = 2
        $cell.Borders.Item(10).LineStyle = 1
    }
}

$lastRow = $startRow + $rows.Count - 1
$ws.Range("A" + ($lastRow + 1)).Select()
